$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 4 (done first, before any inserts shuffle bookmark ordering): drop
# the old "_GoBack" bookmark pair that used to sit right after the json.html
# hyperlink. The bookmark re-appears later, in the new Free-Text Search
# paragraph added below (Change 2).
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

function Wrap-Ooxml($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Change 1: Title paragraph - split "...A0196982H's submission" into 3 runs,
# inserting a new "-A0139744W" run between author id and "'s submission".
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$fullRange = $titlePara.Range
# range excluding the trailing paragraph mark, so the <w:p> attrs/pPr survive
$innerRange = $d.Range($fullRange.Start, $fullRange.End - 1)

$titleRuns = @'
<w:body><w:p><w:r><w:t>This is the README file for A0140713U-A0196982H</w:t></w:r><w:r><w:t>-A0139744W</w:t></w:r><w:r><w:t>&apos;s submission</w:t></w:r></w:p></w:body>
'@
$innerRange.InsertXML((Wrap-Ooxml $titleRuns))

# ---------------------------------------------------------------------------
# Change 2: Replace the empty paragraph right after "Free-Text Search:" with
# two new descriptive paragraphs (and move the _GoBack bookmark here).
# ---------------------------------------------------------------------------
$freeTextSearchPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Free-Text Search:") {
        $freeTextSearchPara = $p
        break
    }
}
$emptyPara = $freeTextSearchPara.Next()
$emptyRange = $emptyPara.Range

$freeTextBody = @'
<w:body>
<w:p>
<w:r><w:t>The free-text search</w:t></w:r>
<w:r><w:t xml:space="preserve"> retrieves documents deemed relevant to the query based on the Vector Space Model and outputs the </w:t></w:r>
<w:r><w:t>all</w:t></w:r>
<w:r><w:t xml:space="preserve"> relevant documents in decreasing order of relevance.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">A posting contains the document ID and the term frequencies in the court name, title, date and content. These </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>tf</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> are handled separately. The scoring scheme implemented for the content and date calculates the score of relevance based on the </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>lnc.ltc</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> score and further awards 1 more point for each successful word match in the query. However, for </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>docID</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>, court and title, a simple score of 2 is given per match in the query. The scores for each of these are stored in a list in the order [</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>docID</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve">, court, title, date, content] as deemed by order of descending important. The lists are then sorted by their natural order to reflect the relevance of the document to the query. Thus, the reason why </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>docID</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve">, court and title are not subjected to </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>tf-idf</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> is because a discrete categorisation is required for the tiered ranking system to work optimally. Only the date and the content are continuous since they need to differentiate between 2 documents.</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
'@
$emptyRange.InsertXML((Wrap-Ooxml $freeTextBody))

# ---------------------------------------------------------------------------
# Change 3: Rework the "Statement of individual work" certify paragraphs to
# add the second author, switch "I" to "we", and drop the gramStart/gramEnd
# proofing markers around "In particular, I".
# ---------------------------------------------------------------------------
$certifyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^\[X\] I, A0140713U-A0196982H, certify") {
        $certifyPara = $p
        break
    }
}
$retrievalPara = $certifyPara.Next()
$expresslyPara = $retrievalPara.Next()
$certifyFullRange = $d.Range($certifyPara.Range.Start, $expresslyPara.Range.End)

$certifyBody = @'
<w:body>
<w:p>
<w:r><w:t>[X] I, A0140713U-A0196982H</w:t></w:r>
<w:r><w:t>-</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>A0139744W</w:t></w:r>
<w:r><w:t xml:space="preserve">, certify that </w:t></w:r>
<w:r><w:t>we</w:t></w:r>
<w:r><w:t xml:space="preserve"> have followed the CS 3245 Information</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">Retrieval class guidelines for homework assignments.  In particular, </w:t></w:r>
<w:r><w:t>we</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">expressly vow that </w:t></w:r>
<w:r><w:t>we</w:t></w:r>
<w:r><w:t xml:space="preserve"> have followed the Facebook rule in discussing</w:t></w:r>
</w:p>
</w:body>
'@
$certifyFullRange.InsertXML((Wrap-Ooxml $certifyBody))
